# Apply the "add cache of rows to improve random access/rewind perfs"
# regression-data edit to data01.xlsx:
#  - Sheet1: B3 becomes 42, new D5 = "bar", selection moves to D6
#  - Sheet2: was empty, gets A1 = "+++Hihi!" (entered with a leading
#    apostrophe so Excel stores it with a quote-prefix style, since the
#    text starts with '+' and would otherwise look like a formula),
#    and C7 = "(+c7+)"; selection moves to C8
#  - Sheet3 / Sheet4: untouched (only cosmetic re-save formatting changes
#    in the original diff, nothing addressable through the object model)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2 first, so the two brand-new shared strings land in the same
#     order as the target file (index 6 = "+++Hihi!", 7 = "(+c7+)") ---
$ws2.Range("A1").Value = "'+++Hihi!"
$ws2.Range("C7").Value = "(+c7+)"

# --- Sheet1 edits ---
$ws1.Range("B3").Value = 42
$ws1.Range("D5").Value = "bar"

# --- Selections: touch each sheet's own selection, then leave Sheet1 as
#     the active/selected tab (matches tabSelected="1" staying on Sheet1) ---
$ws2.Activate() | Out-Null
$ws2.Range("C8").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("D6").Select() | Out-Null
